# Applies the style changes described by the diff against
# test/docx/golden/raw-blocks.docx:
#   1. Add a new paragraph style "AbstractTitle" (based on Normal, next
#      paragraph style "Abstract") positioned ahead of the "Abstract"
#      paragraph it introduces.
#   2. Change the "Abstract" paragraph style's spacing-before from 300
#      to 100 (twentieths of a point -> 15pt -> 5pt).
#   3. Give the "ImportTok" character style a green, bold color.
#   4. Give the "BuiltInTok" character style a green color.

$d = $word.ActiveDocument

# 1. New "Abstract Title" paragraph style.
$abstractTitle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
# Word COM colors are 0x00BBGGRR -- target hex color 345A8A (R=0x34 G=0x5A B=0x8A).
$abstractTitle.Font.Color = 0x34 + (0x5A * 256) + (0x8A * 65536)

# 2. "Abstract" style: spacing before 300 -> 100 (twentieths of a point,
#    i.e. 15pt -> 5pt).
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# 3. "ImportTok" character style: green (008000) + bold.
$importTok = $d.Styles("ImportTok")
$importTok.Font.Color = 0x00 + (0x80 * 256) + (0x00 * 65536)
$importTok.Font.Bold = $true

# 4. "BuiltInTok" character style: green (008000).
$builtInTok = $d.Styles("BuiltInTok")
$builtInTok.Font.Color = 0x00 + (0x80 * 256) + (0x00 * 65536)
